$d = $word.ActiveDocument

function Get-ParaIndexByText($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($needle)) {
            return $i
        }
    }
    return -1
}

function Delete-AlternatePairBeforeContent($contentText) {
    # Deletes the "Alternate:" (Subtitle) paragraph directly preceding the
    # paragraph that starts with $contentText, together with that content
    # paragraph itself - i.e. the whole "Alternate:" / alt-text pair.
    $idx = Get-ParaIndexByText($contentText)
    if ($idx -lt 0) { throw "Could not find paragraph starting with: $contentText" }
    $altIdx = $idx - 1
    $pAlt = $d.Paragraphs($altIdx)
    if (-not $pAlt.Range.Text.StartsWith("Alternate:")) {
        throw "Expected paragraph before '$contentText' to be 'Alternate:'"
    }
    $pContent = $d.Paragraphs($idx)
    $rng = $d.Range($pAlt.Range.Start, $pContent.Range.End)
    $rng.Delete()
}

# 1) Remove the three now-unwanted "Alternate:" + alternate-text paragraph pairs.
#    Go bottom-to-top so earlier paragraph indices stay valid between calls.
Delete-AlternatePairBeforeContent("One too many bad meals at C4C")
Delete-AlternatePairBeforeContent("To increase satisfaction in dining")
Delete-AlternatePairBeforeContent("Is Build Your Own Burrito day")

# 2) Rewrite the "Ever wonder..." paragraph to drop the trailing _GoBack bookmark.
$idx = Get-ParaIndexByText("Ever wonder which wild animals")
$rng = $d.Paragraphs($idx).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Ever wonder which wild animals might be nearby? </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Pokemon</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> aren&#8217;t real but you can pretend that hawk is a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Fearow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> and log him in your personal wildlife journal! Our website will tell you, based on your location, just how safe from nature you are!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml)

# 3) Rewrite "See if you can catch..." paragraph so the leading "S" is its own run.
$idx = Get-ParaIndexByText("See if you can catch")
$rng = $d.Paragraphs($idx).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>S</w:t></w:r><w:r><w:t>ee if you can catch (a glimpse of) &#8216;</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>em</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> all!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml)

# 4) "Motivation:" heading loses its lastRenderedPageBreak.
$idx = Get-ParaIndexByText("Motivation:")
$rng = $d.Paragraphs($idx).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Motivation:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml)

# 5) "Risks:" heading gains the _GoBack bookmark (moved down from the description text).
$idx = Get-ParaIndexByText("Risks:")
$rng = $d.Paragraphs($idx).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Risks:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml)

# 6) "Risk Mitigation Plan:" heading gains the lastRenderedPageBreak that used to sit on "Motivation:".
$idx = Get-ParaIndexByText("Risk Mitigation Plan:")
$rng = $d.Paragraphs($idx).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Risk Mitigation Plan:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml)

# 7) "Proposed Architecture:" heading loses its lastRenderedPageBreak.
$idx = Get-ParaIndexByText("Proposed Architecture:")
$rng = $d.Paragraphs($idx).Range
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml" pkg:padding="512"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Proposed Architecture:</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$rng.InsertXML($xml)

Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
